$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 72: add end-of-day time (HORA F / C72) and PAUSAS (E72); update ASSUNTO/PRODUCAO text
$ws.Range("C72").Value = 0.625
$ws.Range("E72").Value = 0.12222222222222223
$ws.Range("G72").Value = "ESTÁGIO + HARD"
$ws.Range("H72").Value = "Estágio + Hard"

# Row 73: explicit zeros for HORA I / HORA F (previously blank)
$ws.Range("B73").Value = 0
$ws.Range("C73").Value = 0

# Row 74: add start/end-of-day times (HORA I / HORA F) and PAUSAS; update ASSUNTO/PRODUCAO text
$ws.Range("B74").Value = 0.625
$ws.Range("C74").Value = 0.91666666666666663
$ws.Range("E74").Value = 0.092361111111111116
$ws.Range("G74").Value = "HARD"
$ws.Range("H74").Value = "Hard + Atividade voluntária no Alpha EdTech "

# Update the view state to match the edited workbook (scrolled down, selection moved)
$ws.Range("H77").Select()
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
